$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "66.253.46"
Set-TextValue "E2" "  -1.53%  "
Set-TextValue "D3" "3.278.84"
Set-TextValue "E3" "  -1.81%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "580.08"
Set-TextValue "E5" "  +0.33%  "
Set-TextValue "D6" "179.38"
Set-TextValue "E6" "  -2.63%  "
Set-TextValue "D7" "0.629"
Set-TextValue "E7" "  +4.39%  "
Set-TextValue "E8" "  +0.03%  "
Set-TextValue "E9" "  -2.32%  "
Set-TextValue "E10" "  +1.28%  "
Set-TextValue "E11" "  -0.87%  "
Set-TextValue "D12" "3.850.51"
Set-TextValue "E12" "  -1.86%  "
Set-TextValue "E13" "  -3.57%  "
Set-TextValue "D14" "66.277.08"
Set-TextValue "E14" "  -1.76%  "
Set-TextValue "D15" "26.35"
Set-TextValue "E15" "  -3.50%  "
Set-TextValue "E16" "  -1.95%  "
Set-TextValue "D17" "3.291.51"
Set-TextValue "E17" "  -1.61%  "
Set-TextValue "D18" "435.45"
Set-TextValue "E18" "  -2.18%  "
Set-TextValue "D19" "5.52"
Set-TextValue "E19" "  -2.37%  "
Set-TextValue "D20" "13.17"
Set-TextValue "E20" "  -3.04%  "
Set-TextValue "E21" "  -4.25%  "
Set-TextValue "D22" "71.77"
Set-TextValue "E22" "  -2.99%  "
Set-TextValue "E23" "  +0.19%  "
Set-TextValue "D24" "3.422.73"
Set-TextValue "E24" "  -1.97%  "
Set-TextValue "D25" "0.506"
Set-TextValue "E25" "  -1.06%  "
Set-TextValue "D26" "0.198"
Set-TextValue "E26" "  +2.26%  "
Set-TextValue "E27" "  -6.10%  "
Set-TextValue "E28" "  -2.10%  "
Set-TextValue "E29" "  -0.07%  "
Set-TextValue "E30" "  -1.52%  "
Set-TextValue "D31" "22.28"
Set-TextValue "E31" "  -2.70%  "
Set-TextValue "E32" "  +0.10%  "
Set-TextValue "E33" "  -2.47%  "
Set-TextValue "D34" "6.58"
Set-TextValue "E34" "  -2.78%  "
Set-TextValue "E35" "  -3.89%  "
Set-TextValue "D36" "157.66"
Set-TextValue "E36" "  -2.50%  "
Set-TextValue "E37" "  -4.72%  "
Set-TextValue "D38" "26.54"
Set-TextValue "E38" "  -3.58%  "
Set-TextValue "E39" "  -3.13%  "
Set-TextValue "D40" "2.775.19"
Set-TextValue "E40" "  -2.10%  "
Set-TextValue "E41" "  -2.28%  "
Set-TextValue "D42" "4.31"
Set-TextValue "E42" "  -3.36%  "
Set-TextValue "E43" "  -0.13%  "
Set-TextValue "E44" "  -3.31%  "
Set-TextValue "D45" "0.0659"
Set-TextValue "E45" "  -1.83%  "
Set-TextValue "D46" "322.17"
Set-TextValue "E46" "  -0.27%  "
Set-TextValue "D47" "2.29"
Set-TextValue "E47" "  -2.87%  "
Set-TextValue "D48" "23.21"
Set-TextValue "E48" "  -5.28%  "
Set-TextValue "E49" "  -2.51%  "
Set-TextValue "E50" "  +2.97%  "
Set-TextValue "E51" "  +0.04%  "
